$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 29
$ws.Range("B2").Value = 30
$ws.Range("B3").Value = 69
$ws.Range("B4").Value = 96
$ws.Range("B5").Value = 118
$ws.Range("B6").Value = 130
$ws.Range("B7").Value = 148
